$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - reorder recorder emails
$ws.Range("G2").Value = "System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 3 - reorder recorder emails
$ws.Range("G3").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 4 - reorder recorder emails
$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 5 - reorder recorder emails
$ws.Range("G5").Value = "eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 6 - Recorded Sessions count
$ws.Range("L6").Value = 17

# Row 7 - reorder recorder emails
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg"

# Row 8 - Pending Sessions count
$ws.Range("L8").Value = 11

# Row 9 - reorder recorder emails
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# Row 9 - Coverage % (stored as text, leading apostrophe keeps it text instead of a number)
$ws.Range("L9").Value = "'58.6%"

# Row 10 - Average Attendance % (stored as text)
$ws.Range("L10").Value = "'23.4%"

# Row 12 - reorder recorder emails
$ws.Range("G12").Value = "amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

# Row 15 - Group statistics summary row
$ws.Range("O15").Value = 17
$ws.Range("Q15").Value = 11
$ws.Range("R15").Value = "'58.6%"
$ws.Range("S15").Value = "'23.4%"

# Row 16 - session now recorded (was pending) -> copy the "Recorded" (green) look from row 15
$ws.Range("A15:I15").Copy()
$ws.Range("A16:I16").PasteSpecial(-4122)
$ws.Range("G16").Value = "mohamed.saleem@med.asu.edu.eg"
$ws.Range("H16").Value = "22/251"
$ws.Range("I16").Value = "Recorded"

# Row 28 - reorder recorder emails
$ws.Range("G28").Value = "maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
